$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: row, A(col1 value not incl.), B, C, D for rows 5..56
$data = @(
    @(5, 3, 1209.564221739463, 2, 0.08234972050040781),
    @(6, 4, 1209.564221739463, 4, 0.01516576261351482),
    @(7, 5, 1248, 2, 0.07815507531638952),
    @(8, 6, 1266.460117430137, 3, 0.006382493085096246),
    @(9, 7, 1292.473979699731, 2, 0.02049772367939864),
    @(10, 8, 1304.585912295229, 2, 0.02139510558063962),
    @(11, 9, 1317.3, 0, 0.004773763494720345),
    @(12, 10, 1339.018418868032, 2, 0.007785429305680182),
    @(13, 11, 1358.3, 2, 0.4491463597744076),
    @(14, 12, 1452.267230934101, 5, 0.1853376782608937),
    @(15, 13, 1468.498351606071, 4, 0.08373005724241933),
    @(16, 14, 1494.44871075181, 5, 0.1900417373909443),
    @(17, 15, 1515.596556909176, 0, 0.3306340570759647),
    @(18, 16, 1515.596556909176, 2, 0.03880439038249538),
    @(19, 17, 1538.408632604229, 4, 0.04823577340117313),
    @(20, 18, 1574, 4, 0.07289689843887331),
    @(21, 19, 1597, 5, 0.4214536699028364),
    @(22, 20, 1616.872511864008, 3, 0.1079809019393334),
    @(23, 21, 1623.064919830593, 1, 0.002364844158558736),
    @(24, 22, 1642.048348574342, 0, 0.0270023364094122),
    @(25, 23, 1732.134406941665, 4, 0.04641980476414791),
    @(26, 24, 1747.835619913894, 4, 0.05870067363643368),
    @(27, 25, 1772.725731049579, 2, 0.006562677749158311),
    @(28, 26, 1785.161835361124, 2, 0.03823531669500906),
    @(29, 27, 1795.477475116088, 2, 0.01013878774325021),
    @(30, 28, 1805.922029018604, 3, 0.01365473893432243),
    @(31, 29, 1818, 0, 0.08069046492397715),
    @(32, 30, 1840, 2, 0.06897927026451577),
    @(33, 31, 1840, 0, 0.01933433709468267),
    @(34, 32, 1841.298544231364, 2, 0.007096152909158444),
    @(35, 33, 1851.005074483267, 0, 0.005867300457096876),
    @(36, 34, 1865.530752678807, 2, 0.05153847163427385),
    @(37, 35, 1876, 2, 0.04793841377196112),
    @(38, 36, 1885.117065551615, 2, 0.01166998750544201),
    @(39, 37, 1896.361077451411, 2, 0.007999904990412195),
    @(40, 38, 1911.792281651251, 2, 0.01694965968195618),
    @(41, 39, 1924.67520477692, 2, 0.01699344876080397),
    @(42, 40, 1933.550109340654, 2, 0.01755077824842062),
    @(43, 41, 1944.234562301843, 2, 0.01052252851291888),
    @(44, 42, 1959.440632592807, 2, 0.01401285418419295),
    @(45, 43, 1970.376107067389, 2, 0.01625312779948391),
    @(46, 44, 1986.766402625133, 2, 0.01049407157396844),
    @(47, 45, 1995.02814271466, 6, 0.18322299904006),
    @(48, 46, 2014.856993695083, 1, 0.001969364201798255),
    @(49, 47, 2048.833064430492, 2, 0.03176164738981578),
    @(50, 48, 2073.838281975398, 2, 0.04098796786726889),
    @(51, 49, 2092.576632909581, 2, 0.02568909853921173),
    @(52, 50, 2113.24, 2, 0.4736005212905512),
    @(53, 51, 2113.24, 3, 0.1342078440113726),
    @(54, 52, 2133.721396755173, 2, 0.2139754723750275),
    @(55, 53, 2133.721396755173, 3, 0.06640859596183137),
    @(56, 54, 2148.119842376009, 2, 0.01106723146435628),
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

# Apply the same formatting (bold, bordered, centered) used by the existing
# column-A cells to the newly added ones by copying formats only.
$ws.Range("A43").Copy() | Out-Null
$ws.Range("A44:A56").PasteSpecial(-4122) | Out-Null
